$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "90.823.77"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.153.36"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.42"
$ws.Range("E5").Value = "  +9.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "640.76"
$ws.Range("E6").Value = "  +2.99%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  +10.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.367"
$ws.Range("E8").Value = "  -3.07%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.726"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.198"
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.49"
$ws.Range("E13").Value = "  +5.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000251"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.62"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.385.49"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.732.92"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.101.58"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("B19").Value = "SuiNetwork"
$ws.Range("C19").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.77"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.50"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.71"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.66"
$ws.Range("E23").Value = "  +9.41%  "
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.04"
$ws.Range("E25").Value = "  -2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "90.99"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.45"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("E28").Value = "  +0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.72"
$ws.Range("E30").Value = "  +6.94%  "
$ws.Range("E31").Value = "  -4.06%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.203"
$ws.Range("E32").Value = "  +33.53%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.27"
$ws.Range("E33").Value = "  +15.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.90"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "519.58"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.151"
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.22"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +4.30%  "
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.810"
$ws.Range("E40").Value = "  -9.42%  "
$ws.Range("E41").Value = "  +5.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.19"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("E43").Value = "  -3.63%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.41"
$ws.Range("E45").Value = "  +47.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.709"
$ws.Range("E47").Value = "  +13.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "151.37"
$ws.Range("E48").Value = "  +1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.94"
$ws.Range("E49").Value = "  +4.33%  "
$ws.Range("E50").Value = "  +8.47%  "
$ws.Range("E51").Value = "  +4.10%  "
